$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New industry sort order (by descending 行业成交额占比分位数) with updated metrics
$data = @(
    @("农林牧渔", 1.100814819335938, 100, 14.29698950109028, 53.33333206176758),
    @("基础化工", 7.893408298492432, 100, 19.57059103974028, 66.66666412353516),
    @("电力设备及新能源", 12.68211269378662, 100, 30.99728601053357, 100),
    @("交通运输", 1.619480729103088, 91.66666412353516, 11.87767313263281, 66.66666412353516),
    @("钢铁", 1.008443355560303, 91.66666412353516, 12.26769804496032, 80),
    @("纺织服装", 0.6950334310531616, 91.66666412353516, 18.14844804048538, 73.33333587646484),
    @("石油石化", 0.7321617603302002, 91.66666412353516, 15.49144312330321, 60),
    @("建材", 0.9400429129600525, 91.66666412353516, 18.35917181127212, 73.33333587646484),
    @("电力及公用事业", 2.877874374389648, 91.66666412353516, 14.87828036010526, 86.66666412353516),
    @("消费者服务", 0.785746693611145, 83.33333587646484, 18.96184805035591, 60),
    @("综合金融", 0.1826094388961792, 83.33333587646484, 16.25471994876862, 60),
    @("煤炭", 0.8427945971488953, 83.33333587646484, 15.90606394906839, 86.66666412353516),
    @("综合", 0.2067595273256302, 75, 15.68680873124496, 66.66666412353516),
    @("银行", 1.718413472175598, 75, 4.317190465472994, 53.33333206176758),
    @("传媒", 3.283727645874023, 75, 20.59184132133212, 60),
    @("建筑", 1.872642517089844, 75, 17.90679322974182, 46.66666793823242),
    @("医药", 5.495838642120361, 58.33333206176758, 13.3894943180239, 46.66666793823242),
    @("轻工制造", 1.148614883422852, 58.33333206176758, 16.81461450819335, 40),
    @("家电", 1.377744317054749, 50, 12.73197610676289, 40),
    @("汽车", 4.358813762664795, 50, 17.31111774904521, 33.33333206176758),
    @("食品饮料", 1.291661977767944, 50, 9.979829556111133, 53.33333206176758),
    @("房地产", 1.33267617225647, 50, 15.99967069445916, 66.66666412353516),
    @("有色金属", 5.556517601013184, 50, 18.88627188682556, 46.66666793823242),
    @("国防军工", 2.131901979446411, 25, 14.19006614744171, 33.33333206176758),
    @("电子", 15.00424766540527, 25, 18.90659512459271, 26.66666603088379),
    @("商贸零售", 0.904113233089447, 25, 13.28111361242035, 26.66666603088379),
    @("计算机", 8.323662757873535, 25, 20.29446222715908, 40),
    @("机械", 8.012394905090332, 16.66666603088379, 19.86914492620938, 33.33333206176758),
    @("通信", 4.446424007415771, 8.333333015441895, 18.72077462419135, 26.66666603088379),
    @("非银行金融", 2.142033100128174, 8.333333015441895, 5.317189846781717, 26.66666603088379)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
